# Publications template update:
#  - shift the existing header text one column to the right (B..O -> C..P)
#    and put a brand-new "Publication Type" header into B1
#  - rename the old "PDF URL" header (now in P1) to "ISSN/ISBN No"
#  - append four new header columns: Journal Link, UGC Approved, Impact Factor, DOI Link
#  - widen column B and size the new trailing columns to match the new template
#  - re-point the selection at B1 (matches the refreshed template)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the current header row (A1:O1) before overwriting anything.
$headers = @()
for ($c = 1; $c -le 15; $c++) {
    $headers += $ws.Cells.Item(1, $c).Value()
}

# 2) Re-write B1:O1 shifted one column to the right (old B1 -> C1, ... old O1 -> P1),
#    working from the right so we never clobber a value before it's been read.
for ($c = 15; $c -ge 2; $c--) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c - 1]
}

# 3) New "Publication Type" header takes over column B.
$ws.Range("B1").Value = "Publication Type(Journal Paper/Conference Proceedings)"

# 5) Four brand-new trailing headers, copying A1's header look (bold white text
#    on the green fill, thin border, centered) onto each new cell. P1 also gets
#    this treatment since it is a freshly-touched cell beyond the old A1:O1 range.
$ws.Range("A1").Copy()
$ws.Range("P1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) The old "PDF URL" header, now sitting in P1, is renamed.
$ws.Range("P1").Value = "ISSN/ISBN No"

$ws.Range("Q1").Value = "Journal Link"
$ws.Range("R1").Value = "UGC Approved"
$ws.Range("S1").Value = "Impact Factor"
$ws.Range("T1").Value = "DOI Link"

# 6) Column widths for the refreshed layout (columns C..O keep their existing
#    widths; only B is resized and the new trailing columns get explicit widths).
$ws.Columns.Item(2).ColumnWidth = 49.83    # B  -> ~50.6640625
$ws.Columns.Item(16).ColumnWidth = 39.17   # P  -> 40
$ws.Columns.Item(17).ColumnWidth = 39.17   # Q  -> 40
$ws.Columns.Item(18).ColumnWidth = 39.17   # R  -> 40
$ws.Columns.Item(19).ColumnWidth = 20.17   # S  -> 21
$ws.Columns.Item(20).ColumnWidth = 19.83   # T  -> ~20.6640625

# 7) Put the selection/active cell on B1, like the refreshed template.
$ws.Range("B1").Select()
